# "update database for lectures"
#
# Before:  User | Content | Response   (Content holds the A1:L6 12-col data,
#                                        referenced by defined name test__4)
# After:   User | Response | Lecture   (the former "Content" sheet is renamed
#                                        to "Lecture" and moved to the last
#                                        tab; "Response" shifts up to 2nd)

$wb = $excel.ActiveWorkbook

# --- 1. Rename "Content" -> "Lecture" -------------------------------------
$lecture = $wb.Worksheets.Item("Content")
$lecture.Name = "Lecture"

# --- 2. Reorder tabs: User, Response, Lecture -----------------------------
# (move "Lecture" to sit right after "Response")
$response = $wb.Worksheets.Item("Response")
$lecture.Move($null, $response)

# Re-fetch live references by name now that the tab order has changed -
# stale handles captured before a Move/reorder can point at the wrong tab.
$user     = $wb.Worksheets.Item("User")
$response = $wb.Worksheets.Item("Response")
$lecture  = $wb.Worksheets.Item("Lecture")

# --- 3. Restore each sheet's selection / active cell ----------------------
$user.Activate()
$user.Range("G2").Select() | Out-Null

$response.Activate()
$response.Range("B2").Select() | Out-Null

# Lecture (formerly Content) ends up the active / selected tab
$lecture.Activate()
$lecture.Range("B4").Select() | Out-Null
